$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1.xml)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# Existing row 2 (file 236f1989-...) is relabeled with the new run's id/timestamp.
$wsOverview.Range("A2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.md"
$wsOverview.Range("G2").Value = "2016-08-12 15:15:30"
$wsOverview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild hyperlinks on this sheet (the runtime appends rather than updates
# hyperlink display text in place, so clear + recreate them all).
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/710146d34e6e66df6e5b17457cc9f42698cd0998/e2e/1e084702-b98b-41f1-9618-c1d5eeb137cc.md", "", "", "e2e\1e084702-b98b-41f1-9618-c1d5eeb137cc.md")

# New row 3 for the newly generated file.
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md"
$wsOverview.Range("B3").Value = "e2e\47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-12 15:15:30"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/710146d34e6e66df6e5b17457cc9f42698cd0998/e2e/47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md", "", "", "e2e\47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2.xml)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

# Existing row 2 relabeled to the new run's id/timestamps.
$wsZh.Range("A2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.md"
$wsZh.Range("G2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.1c0678f1c576d17149f069c73ac4b0ea93cf3ad6.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-12 15:15:23"
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.md"
$wsZh.Range("J2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.1c0678f1c576d17149f069c73ac4b0ea93cf3ad6.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-12 15:15:50"
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Cells.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/710146d34e6e66df6e5b17457cc9f42698cd0998/e2e/1e084702-b98b-41f1-9618-c1d5eeb137cc.md", "", "", "1e084702-b98b-41f1-9618-c1d5eeb137cc.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0868da1d542354f11a681d6259fbd2bb27c32450/e2e/1e084702-b98b-41f1-9618-c1d5eeb137cc.md", "", "", "1e084702-b98b-41f1-9618-c1d5eeb137cc.md")

# New row 3 for the newly generated file.
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$wsZh.Range("A3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.206f49eb9afc80a7d3c311f31f635da3128db35c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-12 15:15:23"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md"
$wsZh.Range("J3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.206f49eb9afc80a7d3c311f31f635da3128db35c.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-12 15:15:50"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = "'"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = "'"
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/710146d34e6e66df6e5b17457cc9f42698cd0998/e2e/47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md", "", "", "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0868da1d542354f11a681d6259fbd2bb27c32450/e2e/47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md", "", "", "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3.xml)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

# Existing row 2 relabeled to the new run's id/timestamps.
$wsDe.Range("A2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.md"
$wsDe.Range("G2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.1c0678f1c576d17149f069c73ac4b0ea93cf3ad6.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-12 15:15:30"
$wsDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.md"
$wsDe.Range("J2").Value = "1e084702-b98b-41f1-9618-c1d5eeb137cc.1c0678f1c576d17149f069c73ac4b0ea93cf3ad6.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-12 15:16:03"
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Cells.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/710146d34e6e66df6e5b17457cc9f42698cd0998/e2e/1e084702-b98b-41f1-9618-c1d5eeb137cc.md", "", "", "1e084702-b98b-41f1-9618-c1d5eeb137cc.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/67b5709e0e3fd05f194d646ef01fc28b49cf5127/e2e/1e084702-b98b-41f1-9618-c1d5eeb137cc.md", "", "", "1e084702-b98b-41f1-9618-c1d5eeb137cc.md")

# New row 3 for the newly generated file.
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$wsDe.Range("A3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.206f49eb9afc80a7d3c311f31f635da3128db35c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-12 15:15:30"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md"
$wsDe.Range("J3").Value = "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.206f49eb9afc80a7d3c311f31f635da3128db35c.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-12 15:16:03"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = "'"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = "'"
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/710146d34e6e66df6e5b17457cc9f42698cd0998/e2e/47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md", "", "", "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/67b5709e0e3fd05f194d646ef01fc28b49cf5127/e2e/47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md", "", "", "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md")

Write-Host "edit complete"
